$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Rs. 1,499"

$ws.Range("A3").Value = "Truke Buds Q1 In Ear Bluetooth Earphone 8 Hours Playback Bluetooth IPX4(Splash Proof) Powerfull Bass -Bluetooth V 5.1 Black"
$ws.Range("B3").Value = "Rs. 1,199"

$ws.Range("A4").Value = "Syska HE5700 In Ear Bluetooth Neckband 6 Hours Playback IPX4(Splash & Sweat Proof) Powerfull bass -Bluetooth V 5.0 Gray"
$ws.Range("B4").Value = "Rs. 749"

$ws.Range("A5").Value = "NBOX Buzz TWS On Ear True Wireless (TWS) 20 Hours Playback IPX5(Splash & Sweat Proof) Passive noise cancellation -Bluetooth Version 5.1 Black"
$ws.Range("B5").Value = "Rs. 899"
